$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Set header values in the newly inserted row
$ws.Range("A1").Value = "State"
$ws.Range("B1").Value = "change"
$ws.Range("C1").Value = "code2"

# Update selection to A2, matching the target workbook
$ws.Range("A2").Select()
